# Data Science Project Showcase.docx - apply edits
$d = $word.ActiveDocument

# Helper: replace the visible text of a paragraph (excluding its end-of-paragraph
# mark) with new text, collapsing all of its runs into a single run.
function Set-ParagraphText($para, [string]$text) {
    $s = $para.Range.Start
    $e = $para.Range.End - 1
    $rng = $d.Range($s, $e)
    $rng.Text = $text
}

# Helper: insert a brand-new list paragraph right after $afterPara, at list
# level $level (1 = top level / ilvl=0, 2 = sub level / ilvl=1), containing
# $text. Returns the new Paragraph object.
function Insert-ListParagraphAfter($afterPara, [int]$level, [string]$text) {
    $afterPara.Range.InsertParagraphAfter() | Out-Null
    $newPara = $afterPara.Next()
    $newPara.Range.ListFormat.ListLevelNumber = $level
    Set-ParagraphText $newPara $text
    return $newPara
}

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so earlier paragraph indices
# remain valid as the document is edited.
# ---------------------------------------------------------------------------

# 15: "TBD"  -> delete this paragraph entirely (it disappears from the doc)
$p15 = $d.Paragraphs.Item(15)
$p15.Range.Delete()

# 14: "TBD " -> unchanged, leave as-is

# 12/13: "Random Forest Model for Predicting " + "the IMDb Rating of New Movies" (L0)
#        followed by "This model predicts the IMDb rating..." (L1)
#        -> becomes a single "TBD" paragraph (L0); the L1 description paragraph
#           is removed (its text is reused below as the new description for
#           the "Regression Model for Predicting New Movie Ratings" bullet).
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Delete()
$p12 = $d.Paragraphs.Item(12)
Set-ParagraphText $p12 "TBD"

# 11: "F" -> "Regression Model for Predicting New Movie Ratings"
#     + new L1 description paragraph inserted after it.
$p11 = $d.Paragraphs.Item(11)
Set-ParagraphText $p11 "Regression Model for Predicting New Movie Ratings"
Insert-ListParagraphAfter $p11 2 "This model predicts the IMDb rating for new films allowing directors to focus on films predicted to be good ones." | Out-Null

# 10: "E" -> "Application for Parsing Emails in HTML"
#     + new L1 description paragraph inserted after it.
$p10 = $d.Paragraphs.Item(10)
Set-ParagraphText $p10 "Application for Parsing Emails in HTML"
Insert-ListParagraphAfter $p10 2 "This application converted emails in HTML format to readable text then sorted each by email part (from, to, title, etc)." | Out-Null

# 9: "D" -> "Airline Executive Summary Detailing a Campaign"
#    + new L1 description paragraph inserted after it.
$p9 = $d.Paragraphs.Item(9)
Set-ParagraphText $p9 "Airline Executive Summary Detailing a Campaign"
Insert-ListParagraphAfter $p9 2 "This summary presents a proposal to executives of an airline to adopt a PR campaign that would alleviate negative press caused by a crash." | Out-Null

# 8: "C" -> "Image Classifier to Identify Written Numbers"
#    + new L1 description paragraph inserted after it.
$p8 = $d.Paragraphs.Item(8)
Set-ParagraphText $p8 "Image Classifier to Identify Written Numbers"
Insert-ListParagraphAfter $p8 2 "This classifier identifies the names of written numbers as a type of machine-learning-driven OCR (optical character recognition)." | Out-Null

# 7: "This model predicts which deliveries..." -> "This model identifies which deliveries..."
$d.Content.Find.Execute("This model predicts which deliveries", $true, $false, $false, $false, $false, $true, 1, $false, "This model identifies which deliveries", 2) | Out-Null

# 3: weather description paragraph -> insert two new paragraphs after it:
#    new L0 heading "Application for Sorting Biological Data " and
#    new L1 description about the chromatography sorting application.
$p3 = $d.Paragraphs.Item(3)
$newHeading = Insert-ListParagraphAfter $p3 1 "Application for Sorting Biological Data "
Insert-ListParagraphAfter $newHeading 2 "This application sorts data in chromatography files specified by the user and outputs the data as graphs of protein concentration over time. " | Out-Null
